$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial date numbers)
$newRows = @(
    @{ Row = 230; A = 44304; B = 0; C = 2; D = 107.469102632993 },
    @{ Row = 231; A = 44305; B = 0; C = 2; D = 107.469102632993 },
    @{ Row = 232; A = 44306; B = 0; C = 2; D = 107.469102632993 },
    @{ Row = 233; A = 44307; B = 0; C = 2; D = 107.469102632993 }
)

# The preceding row (229) carries the formatting we want to extend
# (column A uses a date-styled cell, columns B-D are unstyled).
$srcRow = 229

foreach ($r in $newRows) {
    $row = $r.Row

    # Only column A carries a cell style (the date format); copy just that
    # cell's formatting down so B-D stay on the default (unstyled) format.
    $ws.Range("A$srcRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}
